# BIS-1002: Fixed XLS export tests
# Adds a new "Internal Assignment" column (O) to the sample-type export
# sheet: a bold header in O4 (matching the other boolean-attribute headers
# K4:N4, but rendered in the larger 12pt header font) and "FALSE" text
# values in O5:O7 (matching the plain data-row style used by the rest of
# the row).

$xlPasteFormats = -4122
$xlPasteValues = -4163

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header cell O4: "Internal Assignment" ------------------------------
# Clone the formatting of the neighbouring boolean-attribute header (K4,
# "Multivalued") - bold Calibri, black - then bump the size to 12pt to
# match the new header font used for this column.
$ws.Range("K4").Copy()
$ws.Range("O4").PasteSpecial($xlPasteFormats)
$ws.Range("O4").Value = "Internal Assignment"
$ws.Range("O4").Font.Size = 12

# --- Data cells O5:O7: "FALSE" ------------------------------------------
# Clone each row's plain data-cell formatting (column A of that row), then
# write a literal text "FALSE" (not a boolean) by writing a text formula
# and converting it to a value in place - this mirrors the plain-text
# "FALSE" entries already used elsewhere in the sheet (e.g. K5:K7).
$rows = @(5, 6, 7)
foreach ($r in $rows) {
    $srcCell = "A" + $r
    $dstCell = "O" + $r

    $ws.Range($srcCell).Copy()
    $ws.Range($dstCell).PasteSpecial($xlPasteFormats)

    $ws.Range($dstCell).Formula = "=""FALSE"""
    $ws.Range($dstCell).Copy()
    $ws.Range($dstCell).PasteSpecial($xlPasteValues)
}

# --- Selection, matching the post-edit active range ----------------------
$ws.Range("O4:O7").Select() | Out-Null
